# Added parser function for lectures:
# - Lecturer sheet: replace "Albert Einstein" with "Igor Miladinovic"
# - Software Engineering sheet: rebuild the lecture table with a new first
#   "Lecture Name" column and rename Start/End Time -> Start/End Date

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("General")
$ws2 = $wb.Worksheets.Item("Lecturer")
$ws4 = $wb.Worksheets.Item("Software Engineering")

# --- Lecturer sheet --------------------------------------------------
$ws2.Range("A2").Value = "Igor Miladinovic"

# --- Software Engineering sheet --------------------------------------
$ws4.Range("H1").Value = "Start Date"
$ws4.Range("I1").Value = "End Date"
